$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 1136, which shifts all the existing rows
# 1136..1227 down to 1137..1228 (matches the diff: dimension grows from
# A1:R1227 to A1:R1228, and every row from 1136 onward carries the data
# that used to belong to the row above it).
$ws.Rows("1136:1136").Insert()

# Populate the newly inserted row 1136 with the new record.
$ws.Range("A1136").Value = 8
$ws.Range("B1136").Value = "Terminal La Palmera de La Serena"
$ws.Range("C1136").Value = "Coquimbo"
$ws.Range("D1136").Value2 = 45223
$ws.Range("E1136").Value = 4
$ws.Range("F1136").Value = 100112004
$ws.Range("G1136").Value = "Cebolla"
$ws.Range("H1136").Value = "Sin especificar"
$ws.Range("I1136").Value = "1a (cosecha)"
$ws.Range("J1136").Value = 2000
$ws.Range("K1136").Value = 12000
$ws.Range("L1136").Value = 13000
$ws.Range("M1136").Value = 12500
$ws.Range("N1136").Value = "$/malla 17 kilos"
$ws.Range("O1136").Value = "Provincia del Elquí"
$ws.Range("P1136").Value = 735
$ws.Range("Q1136").Value = 17
$ws.Range("R1136").Value = "Hortaliza"

# Keep the date-time number format consistent with the rest of column D.
$ws.Range("D1136").NumberFormat = $ws.Range("D1137").NumberFormat
